# Edit applies the diff:
#  1. Rename "Requested quantity" header -> "Weekly_PO_Qty" on "Weekly Quantity" sheet (B1)
#  2. Rename "Requested quantity" header -> "Monthly_PO_Qty" on "Monthly Trend" sheet (B1)
#  3. Add a new "PO Forecast" worksheet (3rd sheet, sheetId 3) with forecast data

$wb = $excel.ActiveWorkbook

$wsWeekly  = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

# --- 1 & 2: header renames ---
$wsWeekly.Range("B1").Value  = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3: add the new "PO Forecast" sheet after the existing sheets ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet  = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "PO Forecast"

# Copy header formatting (bold/bordered style) from the "Weekly Quantity" header row
$wsWeekly.Range("A1:B1").Copy()
$newSheet.Range("A1:D1").PasteSpecial(-4122)

$newSheet.Range("A1").Value = "ds"
$newSheet.Range("B1").Value = "PO_Forecast"
$newSheet.Range("C1").Value = "yhat_lower"
$newSheet.Range("D1").Value = "yhat_upper"

# Copy date-number-format from an existing date cell onto the full A2:A26 date column
$wsWeekly.Range("A2").Copy()
$newSheet.Range("A2:A26").PasteSpecial(-4122)

$newSheet.Range("A2").Value = 44990.99999999999
$newSheet.Range("B2").Value = 75
$newSheet.Range("C2").Value = 18.53915838484487
$newSheet.Range("D2").Value = 125.6161286935401
$newSheet.Range("A3").Value = 44997.99999999999
$newSheet.Range("B3").Value = 72
$newSheet.Range("C3").Value = 19.15340699463327
$newSheet.Range("D3").Value = 126.1263590624558
$newSheet.Range("A4").Value = 45004.99999999999
$newSheet.Range("B4").Value = 70
$newSheet.Range("C4").Value = 10.69650391419025
$newSheet.Range("D4").Value = 121.5921010968618
$newSheet.Range("A5").Value = 45011.99999999999
$newSheet.Range("B5").Value = 67
$newSheet.Range("C5").Value = 14.86204936539071
$newSheet.Range("D5").Value = 118.2570661350112
$newSheet.Range("A6").Value = 45018.99999999999
$newSheet.Range("B6").Value = 65
$newSheet.Range("C6").Value = 13.54290183543577
$newSheet.Range("D6").Value = 117.0573709841987
$newSheet.Range("A7").Value = 45046.99999999999
$newSheet.Range("B7").Value = 54
$newSheet.Range("C7").Value = -0.3029472457527591
$newSheet.Range("D7").Value = 110.2565685592869
$newSheet.Range("A8").Value = 45053.99999999999
$newSheet.Range("B8").Value = 52
$newSheet.Range("C8").Value = -1.242065527296114
$newSheet.Range("D8").Value = 106.3178573921102
$newSheet.Range("A9").Value = 45060.99999999999
$newSheet.Range("B9").Value = 49
$newSheet.Range("C9").Value = -5.925939242430101
$newSheet.Range("D9").Value = 103.4852550359599
$newSheet.Range("A10").Value = 45067.99999999999
$newSheet.Range("B10").Value = 47
$newSheet.Range("C10").Value = -8.022269365892043
$newSheet.Range("D10").Value = 100.6999818994461
$newSheet.Range("A11").Value = 45074.99999999999
$newSheet.Range("B11").Value = 44
$newSheet.Range("C11").Value = -13.73851486262354
$newSheet.Range("D11").Value = 102.1340949701379
$newSheet.Range("A12").Value = 45081.99999999999
$newSheet.Range("B12").Value = 41
$newSheet.Range("C12").Value = -14.37531190089863
$newSheet.Range("D12").Value = 98.35880905888374
$newSheet.Range("A13").Value = 45088.99999999999
$newSheet.Range("B13").Value = 39
$newSheet.Range("C13").Value = -13.06639601860708
$newSheet.Range("D13").Value = 88.20322465324131
$newSheet.Range("A14").Value = 45095.99999999999
$newSheet.Range("B14").Value = 36
$newSheet.Range("C14").Value = -18.54045613929808
$newSheet.Range("D14").Value = 88.05423767033945
$newSheet.Range("A15").Value = 45102.99999999999
$newSheet.Range("B15").Value = 34
$newSheet.Range("C15").Value = -19.76713168991031
$newSheet.Range("D15").Value = 82.72652211958453
$newSheet.Range("A16").Value = 45109.99999999999
$newSheet.Range("B16").Value = 31
$newSheet.Range("C16").Value = -22.16078803734158
$newSheet.Range("D16").Value = 88.07049222474228
$newSheet.Range("A17").Value = 45116.99999999999
$newSheet.Range("B17").Value = 29
$newSheet.Range("C17").Value = -23.43021416621824
$newSheet.Range("D17").Value = 82.37526933365656
$newSheet.Range("A18").Value = 45123.99999999999
$newSheet.Range("B18").Value = 26
$newSheet.Range("C18").Value = -27.75122215814358
$newSheet.Range("D18").Value = 79.81525749741007
$newSheet.Range("A19").Value = 45130.99999999999
$newSheet.Range("B19").Value = 23
$newSheet.Range("C19").Value = -29.14336523117104
$newSheet.Range("D19").Value = 81.06601819640898
$newSheet.Range("A20").Value = 45137.99999999999
$newSheet.Range("B20").Value = 21
$newSheet.Range("C20").Value = -29.85243072476501
$newSheet.Range("D20").Value = 73.48846518251113
$newSheet.Range("A21").Value = 45144.99999999999
$newSheet.Range("B21").Value = 18
$newSheet.Range("C21").Value = -34.85491305673546
$newSheet.Range("D21").Value = 72.5479372664044
$newSheet.Range("A22").Value = 45151.99999999999
$newSheet.Range("B22").Value = 16
$newSheet.Range("C22").Value = -38.12154023283123
$newSheet.Range("D22").Value = 66.95374471411071
$newSheet.Range("A23").Value = 45158.99999999999
$newSheet.Range("B23").Value = 13
$newSheet.Range("C23").Value = -38.54094249405057
$newSheet.Range("D23").Value = 70.47766657876917
$newSheet.Range("A24").Value = 45165.99999999999
$newSheet.Range("B24").Value = 10
$newSheet.Range("C24").Value = -39.84405904504189
$newSheet.Range("D24").Value = 67.3441254584487
$newSheet.Range("A25").Value = 45172.99999999999
$newSheet.Range("B25").Value = 8
$newSheet.Range("C25").Value = -45.87480527460922
$newSheet.Range("D25").Value = 67.06748153924471
$newSheet.Range("A26").Value = 45179.99999999999
$newSheet.Range("B26").Value = 5
$newSheet.Range("C26").Value = -49.69373498513114
$newSheet.Range("D26").Value = 58.06309909536284

Write-Output "PO Forecast sheet populated"
